$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
Write-Output "Rows: $($tbl.Rows.Count)"
Write-Output "Cols: $($tbl.Columns.Count)"
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $tbl.Cell($r, $c)
        $text = $cell.Range.Text
        Write-Output "Cell($r,$c): [$text]"
    }
}
